# Fruta / hortaliza, semanal
# Insert a new weekly record at row 16, shifting the existing rows 16-29 down
# to rows 17-30 (dimension grows from A1:R29 to A1:R30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 16, pushing rows 16..29 to 17..30
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new weekly data point
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 45141
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = "Jengibre"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 400
$ws.Cells.Item(16, 11).Value = 16000
$ws.Cells.Item(16, 12).Value = 17000
$ws.Cells.Item(16, 13).Value = 16550
$ws.Cells.Item(16, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1273
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = "Hortaliza"
